# Update the cryptos list (Price column D, Volume(1h) column E) with the
# latest scraped values, as produced by the scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "29.950.26";  E = "  +0.83%  " }
    3  = @{ D = "1.634.02";   E = "  +1.80%  " }
    4  = @{                   E = "  -0.12%  " }
    5  = @{ D = "214.57";     E = "  +0.84%  " }
    6  = @{                   E = "  +0.29%  " }
    8  = @{ D = "28.73";      E = "  +2.28%  " }
    9  = @{                   E = "  +2.10%  " }
    10 = @{                   E = "  +1.01%  " }
    11 = @{                   E = "  +0.34%  " }
    12 = @{ D = "1.867.54";   E = "  +1.78%  " }
    13 = @{ D = "1.634.15";   E = "  +1.65%  " }
    14 = @{                   E = "  +2.88%  " }
    15 = @{ D = "9.30";       E = "  +18.17%  " }
    16 = @{                   E = "  +2.67%  " }
    17 = @{ D = "29.957.45";  E = "  +0.83%  " }
    18 = @{ D = "64.14" }
    19 = @{ D = "243.24";     E = "  +0.78%  " }
    20 = @{                   E = "  +0.61%  " }
    21 = @{                   E = "  -0.08%  " }
    22 = @{ D = "9.88";       E = "  +5.10%  " }
    23 = @{                   E = "  +2.55%  " }
    24 = @{ D = "2.14";       E = "  +1.17%  " }
    25 = @{ D = "157.71";     E = "  +1.41%  " }
    26 = @{ D = "15.54";      E = "  +0.52%  " }
    27 = @{                   E = "  +1.41%  " }
    28 = @{                   E = "  +2.31%  " }
    29 = @{ D = "1.00";       E = "  -0.14%  " }
    30 = @{                   E = "  +1.23%  " }
    31 = @{                   E = "  +4.22%  " }
    32 = @{                   E = "  +4.26%  " }
    33 = @{                   E = "  -0.35%  " }
    34 = @{ D = "1.424.25" }
    35 = @{                   E = "  +4.64%  " }
    36 = @{                   E = "  +0.06%  " }
    37 = @{                   E = "  -3.29%  " }
    38 = @{                   E = "  -0.12%  " }
    39 = @{ D = "0.0169";     E = "  +0.50%  " }
    40 = @{ D = "75.87";      E = "  +14.42%  " }
    41 = @{ D = "0.551";      E = "  +0.70%  " }
    42 = @{ D = "2.00";       E = "  +2.48%  " }
    43 = @{ D = "0.831";      E = "  +1.62%  " }
    44 = @{                   E = "  -1.36%  " }
    45 = @{ D = "52.97";      E = "  -6.74%  " }
    46 = @{                   E = "  -0.12%  " }
    47 = @{                   E = "  +3.19%  " }
    48 = @{ D = "1.775.72";   E = "  +2.01%  " }
    49 = @{ D = "5.35";       E = "  -0.73%  " }
    50 = @{ D = "0.0₆0113";   E = "  +7.76%  " }
    51 = @{ D = "89.51";      E = "  +3.34%  " }
}

# Price strings that are purely numeric (e.g. "214.57") get auto-coerced to
# numbers by plain Value-assignment; force those specific cells to keep
# being interpreted as text (matching the source feed's string formatting,
# e.g. "9.30" / "1.00" / "2.00" with preserved trailing zeros) by marking
# them as Text before writing.
$numericLookingRows = @(5, 8, 15, 18, 19, 22, 24, 25, 26, 29, 39, 40, 41, 42, 43, 45, 49, 51)

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $dCell = $ws.Range("D$row")
        if ($numericLookingRows -contains $row) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
